$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 fresh rows right above the trailing "fin" marker row (row 107)
# so the marker gets pushed down to row 111.
$ws.Rows("107:110").Insert()

# Domaine column is the same for all four new "systemes" rows.
$ws.Range("A107:A110").Value = "systemes"

# SYS-003 : Resolution de systemes
$ws.Range("B107").Value = "SYS-003"
$ws.Range("C107").Value = "Résolution de systèmes"

# SYS-004 : Inversion de matrice
$ws.Range("B108").Value = "SYS-004"
$ws.Range("C108").Value = "Inversion de matrice"

# SYS-005 et SYS-006 partagent le meme titre de sujet (pivot de Gauss /
# moindres carres), rempli en derniere etape comme dans le classeur original.
$ws.Range("B109").Value = "SYS-005"
$ws.Range("B110").Value = "SYS-006"
$ws.Range("C109:C110").Value = "Régression linéaire par la méthode des moindres carrés"

# Reflect the new bottom-of-sheet selection/scroll position.
$ws.Range("C111").Select()
$excel.ActiveWindow.ScrollRow = 91
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Left = 2820
